$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.346.90"
$ws.Range("E2").Value = "  -0.25%  "
$ws.Range("D3").Value = "'1.848.09"
$ws.Range("E3").Value = "  -0.20%  "
$ws.Range("D4").Value = "'0.9989"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'240.05"
$ws.Range("E6").Value = "  -0.42%  "
$ws.Range("D7").Value = "'0.9992"
$ws.Range("E7").Value = "  -0.14%  "
$ws.Range("D8").Value = "'0.07599"
$ws.Range("E8").Value = "  -1.20%  "
$ws.Range("E9").Value = "  -1.09%  "
$ws.Range("D10").Value = "'24.65"
$ws.Range("E10").Value = "  +0.25%  "
$ws.Range("D11").Value = "'0.07737"
$ws.Range("E11").Value = "  -0.17%  "
$ws.Range("D12").Value = "'5.015"
$ws.Range("E12").Value = "  -0.29%  "
$ws.Range("D13").Value = "'0.6788"
$ws.Range("E13").Value = "  -0.37%  "
$ws.Range("D14").Value = "'0.00001046"
$ws.Range("E14").Value = "  -4.24%  "
$ws.Range("D15").Value = "'82.97"
$ws.Range("E15").Value = "  -0.86%  "
$ws.Range("D16").Value = "'6.121"
$ws.Range("E16").Value = "  -0.59%  "
$ws.Range("D17").Value = "'29.380.46"
$ws.Range("E17").Value = "  -0.28%  "
$ws.Range("D18").Value = "'228.83"
$ws.Range("E18").Value = "  -0.28%  "
$ws.Range("E19").Value = "  -1.21%  "
$ws.Range("D20").Value = "'0.9988"
$ws.Range("E20").Value = "  -0.18%  "
$ws.Range("D21").Value = "'7.457"
$ws.Range("E21").Value = "  +0.20%  "
$ws.Range("D22").Value = "'0.9986"
$ws.Range("E22").Value = "  -0.22%  "
$ws.Range("D23").Value = "'158.51"
$ws.Range("E23").Value = "  +0.84%  "
$ws.Range("E24").Value = "  -0.27%  "
$ws.Range("D25").Value = "'8.428"
$ws.Range("E25").Value = "  +0.44%  "
$ws.Range("D26").Value = "'17.65"
$ws.Range("E26").Value = "  -0.28%  "
$ws.Range("D27").Value = "'1.444"
$ws.Range("E27").Value = "  +9.74%  "
$ws.Range("D28").Value = "'1.468"
$ws.Range("E28").Value = "  +0.10%  "
$ws.Range("D29").Value = "'0.05599"
$ws.Range("E29").Value = "  -2.12%  "
$ws.Range("D30").Value = "'4.101"
$ws.Range("E30").Value = "  -0.65%  "
$ws.Range("D31").Value = "'4.063"
$ws.Range("E31").Value = "  +0.20%  "
$ws.Range("B32").Value = "LidoDAOToken"
$ws.Range("C32").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D32").Value = "'1.829"
$ws.Range("E32").Value = "  -1.08%  "
$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D33").Value = "'1.159"
$ws.Range("E33").Value = "  -0.38%  "
$ws.Range("D34").Value = "'0.6966"
$ws.Range("E34").Value = "  -1.70%  "
$ws.Range("D35").Value = "'2.582"
$ws.Range("E35").Value = "  -0.27%  "
$ws.Range("D36").Value = "'1.232.16"
$ws.Range("E36").Value = "  +0.91%  "
$ws.Range("D37").Value = "'0.01799"
$ws.Range("E37").Value = "  +0.31%  "
$ws.Range("D38").Value = "'2.731"
$ws.Range("E38").Value = "  -1.66%  "
$ws.Range("D39").Value = "'6.379"
$ws.Range("D40").Value = "'0.9045"
$ws.Range("E40").Value = "  -0.60%  "
$ws.Range("D41").Value = "'0.9986"
$ws.Range("E41").Value = "  -0.19%  "
$ws.Range("D42").Value = "'101.26"
$ws.Range("E42").Value = "  -0.34%  "
$ws.Range("D43").Value = "'65.35"
$ws.Range("E43").Value = "  -1.57%  "
$ws.Range("E44").Value = "  +0.61%  "
$ws.Range("D45").Value = "'0.3993"
$ws.Range("E45").Value = "  -0.61%  "
$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").Value = "'1.680"
$ws.Range("E46").Value = "  -0.18%  "
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'8.963"
$ws.Range("E47").Value = "  -0.16%  "
$ws.Range("D48").Value = "'0.1145"
$ws.Range("E48").Value = "  +0.99%  "
$ws.Range("E49").Value = "  -5.17%  "
$ws.Range("E50").Value = "  -0.30%  "
$ws.Range("D51").Value = "'0.4622"
